$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 699
$ws.Range("I5").Value = 100
$ws.Range("J5").Value = 1298
$ws.Range("K5").Value = 100
$ws.Range("L5").Value = 1298
$ws.Range("M5").Value = 15
$ws.Range("N5").Value = -1528

$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").Value = ""

$ws.Range("H62").Value = 3201.0967
$ws.Range("I62").Value = 2910.5
$ws.Range("J62").Value = 3603.4614
$ws.Range("K62").Value = 2910.5
$ws.Range("L62").Value = 3603.4614
$ws.Range("M62").Value = -2286.5
$ws.Range("N62").Value = -4851.4614

$ws.Range("H65").Value = 3201.0967
$ws.Range("I65").Value = 2910.5
$ws.Range("J65").Value = 3603.4614
$ws.Range("K65").Value = 14552.5
$ws.Range("L65").Value = 18017.307
$ws.Range("M65").Value = -11432.5
$ws.Range("N65").Value = -24257.307

$ws.Range("H74").Value = 3299.7368
$ws.Range("J74").Value = 3990
$ws.Range("L74").Value = 3990
$ws.Range("N74").Value = -5862

$ws.Range("H77").Value = 3299.7368
$ws.Range("J77").Value = 3990
$ws.Range("L77").Value = 19950
$ws.Range("N77").Value = -29310

$ws.Range("H87").Value = 44015.5
$ws.Range("J87").Value = 44015.5
$ws.Range("L87").Value = 44015.5
$ws.Range("N87").Value = -46511.5

$ws.Range("H90").Value = 44015.5
$ws.Range("J90").Value = 44015.5
$ws.Range("L90").Value = 132046.5
$ws.Range("N90").Value = -144526.5

$ws.Range("H132").Value = 3021.2666
$ws.Range("I132").Value = 3121.5557
$ws.Range("K132").Value = 9364.667099999999
$ws.Range("M132").Value = -6834.667099999999

$ws.Range("H138").Value = 11767296
$ws.Range("I138").Value = 38462850
$ws.Range("K138").Value = 115388550
$ws.Range("M138").Value = -115383410

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 25001000
$ws.Range("I6").Value = 50000000
$ws.Range("J6").Value = 2000
$ws.Range("K6").Value = 50000000
$ws.Range("L6").Value = 2000
$ws.Range("M6").Value = -49999827
$ws.Range("N6").Value = -2346

$ws.Range("H32").Value = 2563.54
$ws.Range("I32").Value = 2288.5403
$ws.Range("J32").Value = 4403.923
$ws.Range("K32").Value = 2288.5403
$ws.Range("L32").Value = 4403.923
$ws.Range("M32").Value = -2001.5403
$ws.Range("N32").Value = -4977.923

$ws.Range("H74").Value = 26317508
$ws.Range("I74").Value = 30304912
$ws.Range("J74").Value = 638
$ws.Range("K74").Value = 30304912
$ws.Range("L74").Value = 638
$ws.Range("M74").Value = -30304038
$ws.Range("N74").Value = -2386

$ws.Range("H77").Value = 26317508
$ws.Range("I77").Value = 30304912
$ws.Range("J77").Value = 638
$ws.Range("K77").Value = 151524560
$ws.Range("L77").Value = 3190
$ws.Range("M77").Value = -151520192
$ws.Range("N77").Value = -11926

$ws.Range("H102").Value = 1069.2858
$ws.Range("I102").Value = 1080.8334
$ws.Range("J102").Value = 1000
$ws.Range("K102").Value = 1080.8334
$ws.Range("L102").Value = 1000
$ws.Range("M102").Value = 541.1666
$ws.Range("N102").Value = -4244

$ws.Range("H122").Value = 2180.84
$ws.Range("I122").Value = 1675.45
$ws.Range("K122").Value = 5026.35
$ws.Range("M122").Value = -2576.35

$ws.Range("H138").Value = 46886
$ws.Range("J138").Value = 46886
$ws.Range("L138").Value = 46886
$ws.Range("N138").Value = -57166

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3378.4614
$ws.Range("I105").Value = 3380
$ws.Range("J105").Value = 3375
$ws.Range("K105").Value = 3380
$ws.Range("L105").Value = 3375
$ws.Range("M105").Value = -1633
$ws.Range("N105").Value = -6869

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2306.4355
$ws.Range("I31").Value = 1364.9348
$ws.Range("K31").Value = 1364.9348
$ws.Range("M31").Value = -1069.9348

$ws.Range("H34").Value = 2306.4355
$ws.Range("I34").Value = 1364.9348
$ws.Range("K34").Value = 1364.9348
$ws.Range("M34").Value = -1162.9348

$ws.Range("H99").Value = 23813206
$ws.Range("J99").Value = 50003880
$ws.Range("L99").Value = 50003880
$ws.Range("N99").Value = -50006876

$ws.Range("H126").Value = 23813206
$ws.Range("J126").Value = 50003880
$ws.Range("L126").Value = 150011640
$ws.Range("N126").Value = -150016580

$ws.Range("H132").Value = 1863.9546
$ws.Range("I132").Value = 1418.6046
$ws.Range("K132").Value = 4255.8138
$ws.Range("M132").Value = -1725.8138

$ws.Range("H134").Value = 887.6842
$ws.Range("I134").Value = 792.4
$ws.Range("J134").Value = 1999.3334
$ws.Range("K134").Value = 2377.2
$ws.Range("L134").Value = 5998.0002
$ws.Range("M134").Value = 157.8000000000002
$ws.Range("N134").Value = -11068.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 696.64
$ws.Range("J131").Value = 713.73914
$ws.Range("L131").Value = 2141.21742
$ws.Range("N131").Value = -12221.21742

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3864.85
$ws.Range("I80").Value = 3236.5
$ws.Range("J80").Value = 4283.75
$ws.Range("K80").Value = 3236.5
$ws.Range("L80").Value = 4283.75
$ws.Range("M80").Value = -2238.5
$ws.Range("N80").Value = -6279.75

$ws.Range("H83").Value = 3864.85
$ws.Range("I83").Value = 3236.5
$ws.Range("J83").Value = 4283.75
$ws.Range("K83").Value = 16182.5
$ws.Range("L83").Value = 21418.75
$ws.Range("M83").Value = -11190.5
$ws.Range("N83").Value = -31402.75

$ws.Range("H102").Value = 22730112
$ws.Range("I102").Value = 23812452
$ws.Range("K102").Value = 23812452
$ws.Range("M102").Value = -23810830

$ws.Range("H122").Value = 49383816
$ws.Range("I122").Value = 16667386
$ws.Range("K122").Value = 50002158
$ws.Range("M122").Value = -49999708

$ws.Range("H126").Value = 4787.483
$ws.Range("I126").Value = 3995.7222
$ws.Range("J126").Value = 6083.091
$ws.Range("K126").Value = 11987.1666
$ws.Range("L126").Value = 18249.273
$ws.Range("M126").Value = -9517.1666
$ws.Range("N126").Value = -23189.273

$ws.Range("H132").Value = 105078.2
$ws.Range("I132").Value = 76869.5
$ws.Range("J132").Value = 500000
$ws.Range("K132").Value = 230608.5
$ws.Range("L132").Value = 1500000
$ws.Range("M132").Value = -228078.5
$ws.Range("N132").Value = -1505060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 47624348
$ws.Range("I7").Value = 66669572
$ws.Range("K7").Value = 66669572
$ws.Range("M7").Value = -66669460

$ws.Range("H22").Value = 2223.8
$ws.Range("I22").Value = 1359.7778
$ws.Range("J22").Value = 10000
$ws.Range("K22").Value = 1359.7778
$ws.Range("L22").Value = 10000
$ws.Range("M22").Value = -1064.7778
$ws.Range("N22").Value = -10590

$ws.Range("H27").Value = 2223.8
$ws.Range("I27").Value = 1359.7778
$ws.Range("J27").Value = 10000
$ws.Range("K27").Value = 1359.7778
$ws.Range("L27").Value = 10000
$ws.Range("M27").Value = -1252.7778
$ws.Range("N27").Value = -10214

$ws.Range("H40").Value = 3928.3125
$ws.Range("I40").Value = 3670.6667
$ws.Range("J40").Value = 4701.25
$ws.Range("K40").Value = 3670.6667
$ws.Range("L40").Value = 4701.25
$ws.Range("M40").Value = -3534.6667
$ws.Range("N40").Value = -4973.25

$ws.Range("H55").Value = 207.72223
$ws.Range("I55").Value = 216.15384
$ws.Range("K55").Value = 216.15384
$ws.Range("M55").Value = -43.15384

$ws.Range("H126").Value = 47624348
$ws.Range("I126").Value = 66669572
$ws.Range("K126").Value = 200008716
$ws.Range("M126").Value = -200006246

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 168333.33
$ws.Range("J14").Value = 152500
$ws.Range("L14").Value = 152500
$ws.Range("N14").Value = -152836

$ws.Range("H122").Value = 1279.24
$ws.Range("I122").Value = 1295.0416
$ws.Range("K122").Value = 3885.1248
$ws.Range("M122").Value = -1435.1248

$ws.Range("H126").Value = 1206.909
$ws.Range("I126").Value = 1145.3334
$ws.Range("K126").Value = 3436.0002
$ws.Range("M126").Value = -966.0001999999999

$ws.Range("H140").Value = 44900
$ws.Range("J140").Value = 44900
$ws.Range("L140").Value = 44900
$ws.Range("N140").Value = -55260
